$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1,1)
$cell.Range.Text = "60÷3="
$cell = $t.Cell(1,2)
$cell.Range.Text = "50÷7="
$cell = $t.Cell(1,3)
$cell.Range.Text = "93÷8="
$cell = $t.Cell(1,4)
$cell.Range.Text = "56÷8="
$cell = $t.Cell(1,5)
$cell.Range.Text = "26÷8="

$cell = $t.Cell(5,1)
$cell.Range.Text = "69÷7="
$cell = $t.Cell(5,2)
$cell.Range.Text = "62÷3="
$cell = $t.Cell(5,3)
$cell.Range.Text = "52÷3="
$cell = $t.Cell(5,4)
$cell.Range.Text = "75÷9="
$cell = $t.Cell(5,5)
$cell.Range.Text = "38÷6="

$cell = $t.Cell(9,1)
$cell.Range.Text = "33÷5="
$cell = $t.Cell(9,2)
$cell.Range.Text = "81÷6="
$cell = $t.Cell(9,3)
$cell.Range.Text = "55÷3="
$cell = $t.Cell(9,4)
$cell.Range.Text = "93÷9="
$cell = $t.Cell(9,5)
$cell.Range.Text = "60÷5="

$cell = $t.Cell(13,1)
$cell.Range.Text = "96÷7="
$cell = $t.Cell(13,2)
$cell.Range.Text = "46÷5="
$cell = $t.Cell(13,3)
$cell.Range.Text = "25÷3="
$cell = $t.Cell(13,4)
$cell.Range.Text = "56÷2="
$cell = $t.Cell(13,5)
$cell.Range.Text = "60÷5="

$cell = $t.Cell(17,1)
$cell.Range.Text = "30÷6="
$cell = $t.Cell(17,2)
$cell.Range.Text = "30÷2="
$cell = $t.Cell(17,3)
$cell.Range.Text = "88÷4="
$cell = $t.Cell(17,4)
$cell.Range.Text = "82÷9="
$cell = $t.Cell(17,5)
$cell.Range.Text = "16÷8="
